$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.321.31"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.921.54"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("D5").Value = "0.7403"
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("D6").Value = "244.72"
$ws.Range("E6").Value = "  -2.06%  "
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3140"
$ws.Range("E8").Value = "  -2.68%  "
$ws.Range("D9").Value = "27.24"
$ws.Range("E9").Value = "  -2.77%  "
$ws.Range("D10").Value = "0.06984"
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.07995"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "0.7738"
$ws.Range("E12").Value = "  -2.13%  "
$ws.Range("D13").Value = "1.932.66"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").Value = "5.308"
$ws.Range("D15").Value = "91.65"
$ws.Range("E15").Value = "  -3.02%  "
$ws.Range("D16").Value = "30.318.34"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "14.24"
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("D18").Value = "246.02"
$ws.Range("E18").Value = "  -2.85%  "
$ws.Range("D19").Value = "5.848"
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007858"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.170.49"
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("D23").Value = "1.004"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.670"
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.420"
$ws.Range("E25").Value = "  -1.74%  "
$ws.Range("D26").Value = "165.45"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("D27").Value = "18.98"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("D28").Value = "0.1273"
$ws.Range("E28").Value = "  -5.25%  "
$ws.Range("D29").Value = "2.135"
$ws.Range("E29").Value = "  -8.10%  "
$ws.Range("D30").Value = "1.355"
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").Value = "1.552"
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("D32").Value = "4.352"
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("D33").Value = "4.073"
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("D34").Value = "0.05193"
$ws.Range("E34").Value = "  +1.50%  "
$ws.Range("D35").Value = "1.301"
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").Value = "0.7488"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Value = "2.778"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("E38").Value = "  -1.41%  "
$ws.Range("D39").Value = "2.797"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").Value = "6.358"
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("D41").Value = "75.83"
$ws.Range("E41").Value = "  -3.17%  "
$ws.Range("D42").Value = "0.4466"
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("D43").Value = "1.947"
$ws.Range("E43").Value = "  -2.53%  "
$ws.Range("D44").Value = "1.003"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").Value = "0.8362"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("D46").Value = "7.655"
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "9.925"
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "101.33"
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1220"
$ws.Range("E50").Value = "  +4.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "940.40"
$ws.Range("E51").Value = "  -6.29%  "
